$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 772, shifting rows 772-813 down to 773-814
$ws.Rows.Item(772).Insert()

# Populate the newly inserted row 772 with the new data point
$ws.Cells.Item(772, 1).NumberFormat = "@"
$ws.Cells.Item(772, 1).Value = "2026/02/08"
$ws.Cells.Item(772, 1).Style = "Normal"
$ws.Cells.Item(772, 2).Value = "日"
$ws.Cells.Item(772, 3).Value = 8
$ws.Cells.Item(772, 4).Value = 87
